$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Royal Dutch" product row (original row 2). This shifts the
# "Παπαγάλος" coffee row up to row 2 and the totals row up to row 3.
$ws.Rows(2).Delete()

# Update the (now) single product row with its new campaign dates,
# quantity and turnover.
$ws.Range("C2").Value = 43988
$ws.Range("D2").Value = 43997
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 10.15

# Update the totals row to match the single remaining product row.
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 10.15

# Shrink the conditional-formatting ranges so they stop covering the
# now-removed third row.
$iFc = $ws.Range("I1:I3").FormatConditions.Item(1)
$iFc.ModifyAppliesToRange($ws.Range("I1:I2"))

$jFc = $ws.Range("J1:J3").FormatConditions.Item(1)
$jFc.ModifyAppliesToRange($ws.Range("J1:J2"))
